$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (column F) values for the affected rows
$ws.Range("F5").Value  = 5
$ws.Range("F9").Value  = -1
$ws.Range("F15").Value = -2
$ws.Range("F17").Value = 2
$ws.Range("F18").Value = -1
$ws.Range("F23").Value = -2
$ws.Range("F24").Value = -1
$ws.Range("F31").Value = 2
$ws.Range("F43").Value = 6
$ws.Range("F44").Value = -5
$ws.Range("F45").Value = -1
$ws.Range("F46").Value = -4
